$d = $word.ActiveDocument

# Locate the unique sentence that needs the duration digit changed.
$searchRange = $d.Content
$searchRange.Find.Execute("atnaujintas api (1 val.)", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

$sentenceStart = $searchRange.Start
# "atnaujintas api (" is 17 characters long, so the duration digit "1" sits
# right after it.
$digitStart = $sentenceStart + 17
$digitRange = $d.Range($digitStart, $digitStart + 1)

# Use tracked changes so the replacement keeps its own run instead of being
# silently re-merged into the surrounding text when we accept it -- this
# mirrors how Word naturally ends up with separate runs after an edit.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$digitRange.Delete()
$insertPoint = $d.Range($digitStart, $digitStart)
$insertPoint.InsertAfter("2")

$d.AcceptAllRevisions()
$d.TrackRevisions = $wasTracking
